$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.892.01'
$ws.Range("E2").Value = '  -0.75%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.153.88'
$ws.Range("E3").Value = '  -0.98%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '584.95'
$ws.Range("E5").Value = '  +1.75%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.83'
$ws.Range("E6").Value = '  -2.99%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.155.73'
$ws.Range("E8").Value = '  -0.88%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.528'
$ws.Range("E9").Value = '  -0.16%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.159'
$ws.Range("E10").Value = '  -3.37%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.17'
$ws.Range("E11").Value = '  -1.30%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.500'
$ws.Range("E12").Value = '  -1.70%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000267'
$ws.Range("E13").Value = '  -3.03%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.37'
$ws.Range("E14").Value = '  -2.97%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.667.77'
$ws.Range("E15").Value = '  -0.91%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.894.36'
$ws.Range("E16").Value = '  -0.71%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.16'
$ws.Range("E17").Value = '  -1.36%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.145.71'
$ws.Range("E18").Value = '  -1.07%  '
$ws.Range("E19").Value = '  -0.52%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '503.64'
$ws.Range("E20").Value = '  -2.47%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.54'
$ws.Range("E21").Value = '  +3.60%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.717'
$ws.Range("E22").Value = '  -3.49%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '15.04'
$ws.Range("E23").Value = '  -6.40%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.82'
$ws.Range("E24").Value = '  -1.35%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.52'
$ws.Range("E25").Value = '  -0.74%  '
$ws.Range("E26").Value = '  +0.09%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.09'
$ws.Range("E27").Value = '  -1.28%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.93'
$ws.Range("E28").Value = '  +0.49%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.20'
$ws.Range("E29").Value = '  -1.00%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.83'
$ws.Range("E30").Value = '  -0.11%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '27.71'
$ws.Range("E31").Value = '  -1.70%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.23'
$ws.Range("E32").Value = '  -0.56%  '
$ws.Range("E33").Value = '  +0.02%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.45'
$ws.Range("E34").Value = '  +1.83%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.48'
$ws.Range("E35").Value = '  -2.96%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '54.99'
$ws.Range("E36").Value = '  -1.78%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0898'
$ws.Range("E37").Value = '  +1.88%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '469.72'
$ws.Range("E38").Value = '  -2.30%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0418'
$ws.Range("E39").Value = '  -1.23%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.91'
$ws.Range("E40").Value = '  -7.29%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.77'
$ws.Range("E41").Value = '  +0.90%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.993.50'
$ws.Range("E42").Value = '  -4.66%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.117'
$ws.Range("E43").Value = '  -3.77%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.284'
$ws.Range("E44").Value = '  -2.74%  '
$ws.Range("E45").Value = '  -4.33%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '28.36'
$ws.Range("E46").Value = '  -3.74%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0₃0603'
$ws.Range("E47").Value = '  +2.64%  '
$ws.Range("E48").Value = '  +0.01%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.115'
$ws.Range("E49").Value = '  -1.50%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.25'
$ws.Range("E50").Value = '  -3.71%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '119.46'
$ws.Range("E51").Value = '  -4.17%  '
